# Apply the data refresh captured in the commit "Update gh-pages to output
# generated at 456a3b4": a handful of "want-to-go" counts (column F) were
# bumped, and two events' "lowest ticket price" cells (column G) were
# switched from a numeric price to the text "不可售" (not for sale).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value  = 922
$ws1.Range("G3").Value  = "不可售"
$ws1.Range("F4").Value  = 545
$ws1.Range("F6").Value  = 12
$ws1.Range("F7").Value  = 704
$ws1.Range("F8").Value  = 316
$ws1.Range("F10").Value = 114
$ws1.Range("F11").Value = 230
$ws1.Range("F12").Value = 184
$ws1.Range("F13").Value = 4335
$ws1.Range("F14").Value = 30
$ws1.Range("F15").Value = 21
$ws1.Range("F16").Value = 459
$ws1.Range("F18").Value = 520
$ws1.Range("F19").Value = 309
$ws1.Range("F23").Value = 689
$ws1.Range("F25").Value = 282
$ws1.Range("F26").Value = 981
$ws1.Range("F28").Value = 1653
$ws1.Range("F29").Value = 387

# ---------------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("F5").Value = 251
$ws2.Range("F7").Value = 250

# ---------------------------------------------------------------------
# Sheet "本地生活" (Local Life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$ws3.Range("F3").Value = 143

# ---------------------------------------------------------------------
# Sheet "全部类型" (All types - aggregate of the sheets above)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value  = 922
$ws4.Range("G4").Value  = "不可售"
$ws4.Range("F7").Value  = 143
$ws4.Range("F8").Value  = 545
$ws4.Range("F10").Value = 12
$ws4.Range("F11").Value = 704
$ws4.Range("F13").Value = 316
$ws4.Range("F15").Value = 114
$ws4.Range("F16").Value = 230
$ws4.Range("F17").Value = 184
$ws4.Range("F18").Value = 184
$ws4.Range("F19").Value = 4335
$ws4.Range("F20").Value = 30
$ws4.Range("F21").Value = 21
$ws4.Range("F22").Value = 251
$ws4.Range("F23").Value = 459
$ws4.Range("F25").Value = 520
$ws4.Range("F26").Value = 309
$ws4.Range("F31").Value = 250
$ws4.Range("F34").Value = 689
$ws4.Range("F39").Value = 282
$ws4.Range("F40").Value = 981
$ws4.Range("F42").Value = 1653
$ws4.Range("F43").Value = 387
